$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet,
#    mirroring the layout of the most recent quarter sheet ("2021-Q4").
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Match the page margins used by the other quarterly sheets.
$ps = $newSheet.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# Copy the header/body cell formatting (borders, bold, alignment) from the
# template sheet so the new sheet matches the existing quarterly sheets.
$template.Range("B1:H5").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$template.Range("A2:A5").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Fund holdings data for 2022-Q1
$rows = @(
    @(0, "001302", "前海开源金银珠宝主题精选混合A",       "8.61", "91.91", "7.98", "0.6871", 7),
    @(1, "003304", "前海开源沪港深核心资源灵活配置混合A", "5.91", "93.10", "6.66", "0.3936", 10),
    @(2, "002207", "前海开源金银珠宝主题精选混合C",       "3.45", "91.91", "7.98", "0.2753", 7),
    @(3, "003305", "前海开源沪港深核心资源灵活配置混合C", "2.19", "93.10", "6.66", "0.1459", 10)
)

# Columns B (fund code) and D:G (scale/position/value figures) are stored
# as text in the source data (e.g. to preserve leading zeros), so force a
# text number format before writing them.
$newSheet.Range("B2:B5").NumberFormat = "@"
$newSheet.Range("D2:G5").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $newSheet.Cells.Item($r,1).Value = $row[0]
    $newSheet.Cells.Item($r,2).Value = $row[1]
    $newSheet.Cells.Item($r,3).Value = $row[2]
    $newSheet.Cells.Item($r,4).Value = $row[3]
    $newSheet.Cells.Item($r,5).Value = $row[4]
    $newSheet.Cells.Item($r,6).Value = $row[5]
    $newSheet.Cells.Item($r,7).Value = $row[6]
    $newSheet.Cells.Item($r,8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2. Add a new "2022-Q1" summary row at the top of the "总计" sheet,
#    shifting the existing rows down by one.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")

# Capture the existing summary rows (currently rows 2-6) before rewriting.
$existing = @()
for ($r = 2; $r -le 6; $r++) {
    $existing += ,@($ws.Cells.Item($r,2).Value(), $ws.Cells.Item($r,3).Value(), $ws.Cells.Item($r,4).Value())
}

# Preserve the "index" column (A) styling by copying it down to the new
# last row before we rewrite all the values.
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)

for ($i = 0; $i -lt $existing.Length; $i++) {
    $destRow = $i + 3
    $data = $existing[$i]
    $ws.Cells.Item($destRow, 1).Value = $i + 1
    $ws.Cells.Item($destRow, 2).Value = $data[0]
    $ws.Cells.Item($destRow, 3).Value = $data[1]
    $ws.Cells.Item($destRow, 4).Value = $data[2]
}

# Write the new first summary row for 2022-Q1.
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "2022-Q1"
$ws.Cells.Item(2,3).Value = 4
$ws.Cells.Item(2,4).Value = 1.5
